$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-6 for changed numeric cells
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1618313333333333
$ws.Range("H2").Value = 0.485494
$ws.Range("M2").Value = 6.134954333333333
$ws.Range("N2").Value = 18.404863
$ws.Range("O2").Value = 0.2326051299917221
$ws.Range("P2").Value = 0.232605129991722
$ws.Range("Q2").Value = 0.9928278397024444
$ws.Range("R2").Value = 8.935450557322
$ws.Range("S2").Value = 0.2326051299917221
$ws.Range("T2").Value = 0.232605129991722

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1618313333333333
$ws.Range("H3").Value = 0.485494
$ws.Range("O3").Value = 0.1862453502798354
$ws.Range("P3").Value = 0.1862453502798353
$ws.Range("Q3").Value = 0.794950518845111
$ws.Range("R3").Value = 7.154554669606
$ws.Range("S3").Value = 0.1862453502798354
$ws.Range("T3").Value = 0.1862453502798353

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1618313333333333
$ws.Range("H4").Value = 0.485494
$ws.Range("M4").Value = 0.008733333333333334
$ws.Range("N4").Value = 0.0262
$ws.Range("O4").Value = 0.0003311219651992584
$ws.Range("P4").Value = 0.0003311219651992584
$ws.Range("Q4").Value = 0.001413326977777778
$ws.Range("R4").Value = 0.0127199428
$ws.Range("S4").Value = 0.0003311219651992584
$ws.Range("T4").Value = 0.0003311219651992584

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1618313333333333
$ws.Range("H5").Value = 0.485494
$ws.Range("M5").Value = 14.59271166666667
$ws.Range("N5").Value = 43.778135
$ws.Range("O5").Value = 0.5532787058762763
$ws.Range("P5").Value = 0.5532787058762761
$ws.Range("Q5").Value = 2.361557985965555
$ws.Range("R5").Value = 21.25402187369
$ws.Range("S5").Value = 0.5532787058762763
$ws.Range("T5").Value = 0.5532787058762761

$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1618313333333333
$ws.Range("H6").Value = 0.485494
$ws.Range("M6").Value = 0.5433020000000001
$ws.Range("N6").Value = 1.629906
$ws.Range("O6").Value = 0.02059914800801765
$ws.Range("P6").Value = 0.02059914800801765
$ws.Range("Q6").Value = 0.08792328706266668
$ws.Range("R6").Value = 0.791309583564
$ws.Range("S6").Value = 0.02059914800801765
$ws.Range("T6").Value = 0.02059914800801765

# Add new row 7 (Resolving-Mac)
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Edn3"
$ws.Range("C7").Value = "Ednra"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1618313333333333
$ws.Range("H7").Value = 0.485494
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.1830566666666666
$ws.Range("N7").Value = 0.5491699999999999
$ws.Range("O7").Value = 0.006940543878949493
$ws.Range("P7").Value = 0.006940543878949492
$ws.Range("Q7").Value = 0.02962430444222222
$ws.Range("R7").Value = 0.26661873998
$ws.Range("S7").Value = 0.006940543878949493
$ws.Range("T7").Value = 0.006940543878949492
